$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Yes"
$ws.Range("F3").Value = "jordana"
$ws.Range("F4").Value = "carlos"
$ws.Range("F5").Value = "jordana"
$ws.Range("B6").Value = "No"
$ws.Range("F6").Value = "carlos"

$ws.Range("B5").Select()
